# "Common: Improved UI, a lot"
#
# Appends 29 new Czech ("cs") translation rows to the end of the Import
# sheet - new labels/translations for the lab.mixture.*, lab.build.*,
# lab.vape.*, lab.liquid.* and lab.coil.* admin-UI sections - and leaves
# the selection where Excel would after keying them in (cell B798).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastOldRow = 777
$firstNewRow = 778
$lastNewRow = 806

# Carry the "import" cell style (wrap-text) used throughout the table down
# across the whole new block before writing any values into it, so the new
# rows pick up the same formatting as every other data row.
$ws.Range("A777:C777").Copy()
$ws.Range("A778:C806").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row / column / value triples, in the exact order the cells were populated
# (language column first, then the keys and their translations) so that
# repeated translations - "Seznam", "Novy", "Editace", "Klon", etc. - line
# up with how they were first introduced while typing the table.
$rows = @(778, 779, 780, 781, 782, 783, 784, 785, 786, 787, 788, 789, 790, 791, 792, 793, 794, 795, 796, 797, 798, 799, 800, 801, 802, 803, 804, 805, 806, 778, 779, 784, 788, 793, 796, 799, 801, 804, 778, 779, 780, 781, 782, 783, 782, 789, 795, 800, 805, 783, 780, 785, 791, 797, 802, 781, 786, 792, 798, 784, 785, 786, 787, 787, 794, 806, 788, 789, 790, 790, 791, 792, 793, 794, 795, 796, 797, 798, 800, 799, 801, 802, 803, 803, 804, 805, 806)
$cols = @("A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "A", "C", "C", "C", "C", "C", "C", "C", "C", "C", "B", "B", "B", "B", "B", "B", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "B", "B", "B", "B", "C", "C", "C", "B", "B", "B", "C", "B", "B", "B", "B", "B", "B", "B", "B", "B", "B", "B", "B", "B", "C", "B", "B", "B")
$vals = @("cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "cs", "Výrobce", "Mixy", "Buildy", "Build", "Vape", "Liquidy", "Liquid", "Spirálky", "Spirálka", "lab.mixture.preview.vendor", "lab.mixture.label", "lab.mixture.list.label", "lab.mixture.create.label", "lab.mixture.edit.label", "lab.mixture.index.label", "Editace", "Editace", "Editace", "Editace", "Editace", "Detail", "Seznam", "Seznam", "Seznam", "Seznam", "Seznam", "Nový", "Nový", "Nový", "Nový", "lab.build.label", "lab.build.list.label", "lab.build.create.label", "lab.build.clone.label", "Klon", "Klon", "Klon", "lab.build.index.label", "lab.build.edit.label", "lab.vape.label", "Vapy", "lab.vape.list.label", "lab.vape.create.label", "lab.vape.index.label", "lab.vape.clone.label", "lab.vape.edit.label", "lab.liquid.label", "lab.liquid.list.label", "lab.liquid.create.label", "lab.liquid.edit.label", "lab.liquid.index.label", "lab.coil.label", "lab.coil.list.label", "lab.coil.create.label", "Nová", "lab.coil.index.label", "lab.coil.edit.label", "lab.coil.clone.label")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Range($cols[$i] + $rows[$i]).Value = $vals[$i]
}

# Mirror the cursor position Excel would leave behind after this entry.
$ws.Range("B798").Select()
